$wb = $excel.ActiveWorkbook

$sheetNaidu = $wb.Worksheets.Item("NAIDU")
$sheetMale = $wb.Worksheets.Item("NAIDU-Male")
$sheetFirst371 = $wb.Worksheets.Item("0_First_371_Last_371")
$sheetFirst590 = $wb.Worksheets.Item("1_First_590_Last_590")

# Update Age (column T) from 27 to 28 for the two matching people
$sheetNaidu.Range("T6").Value = 28
$sheetNaidu.Range("T9").Value = 28

$sheetMale.Range("T5").Value = 28
$sheetMale.Range("T8").Value = 28

# Set Status (column U) on row 2 to "NA"
$sheetFirst371.Range("U2").Value = "NA"
$sheetFirst590.Range("U2").Value = "NA"
